$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.761.28'
$ws.Range('E2').Value = '  -1.53%  '
$ws.Range('D3').Value = '2.180.36'
$ws.Range('E3').Value = '  -2.95%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '293.11'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '86.23'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.81%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.557'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.40%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.471'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -9.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '31.80'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -8.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0753'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -7.01%  '
$ws.Range('E12').Value = '  -2.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.63'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -7.15%  '
$ws.Range('D14').Value = '2.516.58'
$ws.Range('E14').Value = '  -2.81%  '
$ws.Range('D15').Value = '2.254.81'
$ws.Range('E15').Value = '  +0.46%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '12.78'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.88%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.752'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -10.17%  '
$ws.Range('D18').Value = '43.335.99'
$ws.Range('E18').Value = '  -1.83%  '
$ws.Range('D19').Value = '0.0₃0865'
$ws.Range('E19').Value = '  -10.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.74'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -9.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.58'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -14.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '62.16'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '225.30'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.74'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -6.97%  '
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.78'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -10.44%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.22'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.04'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.91%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '34.50'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -10.92%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '18.79'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.38%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '145.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.17'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -12.87%  '
$ws.Range('E33').Value = '  -6.96%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0709'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -11.20%  '
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.114'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.37%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.84'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -8.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0994'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.81%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.60'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -11.31%  '
$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.00'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.44'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -10.19%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0275'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -8.84%  '
$ws.Range('B42').Value = 'NEARProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.99'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -13.66%  '
$ws.Range('B43').Value = 'Celestia'
$ws.Range('C43').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.88'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -12.19%  '
$ws.Range('D44').Value = '1.743.90'
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.61'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.20%  '
$ws.Range('B46').Value = 'BitcoinSV'
$ws.Range('C46').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '71.75'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -11.01%  '
$ws.Range('B47').Value = 'ordi'
$ws.Range('C47').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '65.65'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.66%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.167'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -12.90%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.399.51'
$ws.Range('E49').Value = '  -2.78%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '89.90'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -9.63%  '
$ws.Range('B51').Value = 'HuobiToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.66'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +7.68%  '
